# Update the "Customer total sales" sheet so that rows 4 and 5 swap their
# id/name values (WALA771012HCRGR054 / Wednesday Addams moves up to row 4,
# GOTW771012HMRGR087 / Khal Drogo moves down to row 5), and the purchases
# value that used to belong to row 5 (4) becomes 0, since the person now
# in row 4 already had 0 purchases.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A4").Value = "WALA771012HCRGR054"
$ws.Range("B4").Value = "Wednesday Addams"
$ws.Range("C4").Value = 0

$ws.Range("A5").Value = "GOTW771012HMRGR087"
$ws.Range("B5").Value = "Khal Drogo"
$ws.Range("C5").Value = 0
